$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 33 with data matching the existing table pattern
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"

# Update the view to match the new scroll/selection position
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("C31").Select()
